$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.795.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.893.13"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.56"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06552"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07775"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.53"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.890.85"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.38%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7395"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.257"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.87%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.58"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.862.08"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.24"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007597"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.139.04"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.351"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.268"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.277"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.62%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09762"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.27%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.502"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.325"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.204"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04904"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7020"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01916"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.371"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "75.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.035"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4286"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8388"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.09"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.433"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.86"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "928.83"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05765"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.13%  "
